$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.755.72'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.306.84'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.47%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '185.13'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.53%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.300.45'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.47%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.574'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.174'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.575'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.68'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000260'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.839.78'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.40'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '571.56'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -10.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.695.06'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("E18").Value = '  +0.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.311.07'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.79'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.887'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.89'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.97'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '98.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.40'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.64'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.66'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -9.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '555.76'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '10.79'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.51%  '
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.757.57'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.21%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.102'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.58'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '33.48'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.25%  '
$ws.Range("E40").Value = '  -3.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0₃0679'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.11'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.56'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.25%  '
$ws.Range("E44").Value = '  +1.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.330'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0406'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.96'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -14.21%  '
$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.999'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.09%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.126'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.50'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '124.35'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.26%  '
